# Updates cryptos list values (price/volume columns) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.579.33"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "3.509.65"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'586.74"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("D6").Value = "'132.83"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("D7").Value = "3.509.84"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "'0.389"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "4.108.71"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "3.511.08"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").Value = "64.605.59"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'9.97"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "'14.29"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").Value = "'390.07"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").Value = "3.652.63"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").Value = "'74.11"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D27").Value = "'0.0000111"
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("D28").Value = "'1.57"
$ws.Range("E28").Value = "  -4.95%  "
$ws.Range("D29").Value = "'7.46"
$ws.Range("E29").Value = "  -7.59%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").Value = "'8.25"
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("D33").Value = "3.516.22"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "'5.30"
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("D38").Value = "'1.58"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "'171.32"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "'6.98"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "'0.0811"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").Value = "'0.813"
$ws.Range("E42").Value = "  -2.10%  "
$ws.Range("D43").Value = "'26.04"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'42.15"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").Value = "2.467.00"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0269"
$ws.Range("E51").Value = "  +0.12%  "
